# Apply the "Updated cryptos list" data refresh described by the commit diff.
# Each target cell already stores plain text (inline strings) in the workbook,
# so every write below keeps that text semantics. Cells whose new value would
# otherwise be auto-parsed by Excel as a plain number (e.g. "1.00", "5.60")
# get a leading single-quote text-prefix so the stored value keeps its exact
# original formatting (trailing zeros, etc.) instead of being normalized to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.615.96"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").Value = "2.442.77"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "'566.92"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("D6").Value = "'145.92"
$ws.Range("E6").Value = "  +1.73%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.57%  "

# Row 9
$ws.Range("E9").Value = "  +2.49%  "

# Row 10
$ws.Range("E10").Value = "  +0.24%  "

# Row 11
$ws.Range("E11").Value = "  -1.17%  "

# Row 12
$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000185"
$ws.Range("E13").Value = "  +5.54%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'26.87"
$ws.Range("E14").Value = "  +4.27%  "

# Row 15
$ws.Range("D15").Value = "2.820.01"

# Row 16
$ws.Range("D16").Value = "62.214.46"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range("D17").Value = "2.436.92"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18
$ws.Range("D18").Value = "'11.31"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19
$ws.Range("E19").Value = "  +1.79%  "

# Row 20
$ws.Range("D20").Value = "'325.58"
$ws.Range("E20").Value = "  +0.33%  "

# Row 21
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  -0.21%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").Value = "'67.42"
$ws.Range("E23").Value = "  +2.60%  "

# Row 24
$ws.Range("D24").Value = "'1.77"
$ws.Range("E24").Value = "  +2.88%  "

# Row 25
$ws.Range("E25").Value = "  -3.11%  "

# Row 26
$ws.Range("D26").Value = "'568.67"
$ws.Range("E26").Value = "  -1.35%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0985"
$ws.Range("E27").Value = "  +3.46%  "

# Row 28
$ws.Range("D28").Value = "2.562.04"
$ws.Range("E28").Value = "  +1.28%  "

# Row 29
$ws.Range("E29").Value = "  -0.33%  "

# Row 30
$ws.Range("D30").Value = "'8.41"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31
$ws.Range("D31").Value = "'1.46"
$ws.Range("E31").Value = "  +1.19%  "

# Row 32
$ws.Range("E32").Value = "  -0.52%  "

# Row 33
$ws.Range("E33").Value = "  +0.61%  "

# Row 34
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("E35").Value = "  +3.48%  "

# Row 36
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37
$ws.Range("E37").Value = "  +0.37%  "

# Row 38
$ws.Range("D38").Value = "'5.60"
$ws.Range("E38").Value = "  +0.38%  "

# Row 39
$ws.Range("D39").Value = "'18.80"
$ws.Range("E39").Value = "  +0.45%  "

# Row 40
$ws.Range("D40").Value = "'150.50"
$ws.Range("E40").Value = "  -1.64%  "

# Row 41
$ws.Range("E41").Value = "  +0.46%  "

# Row 42
$ws.Range("E42").Value = "  +0.69%  "

# Row 43
$ws.Range("D43").Value = "'2.43"
$ws.Range("E43").Value = "  +5.90%  "

# Row 44
$ws.Range("D44").Value = "'148.91"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("E45").Value = "  +1.34%  "

# Row 46
$ws.Range("E46").Value = "  +0.21%  "

# Row 47
$ws.Range("D47").Value = "'20.46"
$ws.Range("E47").Value = "  +1.50%  "

# Row 48
$ws.Range("E48").Value = "  +1.07%  "

# Row 49
$ws.Range("D49").Value = "'0.0931"
$ws.Range("E49").Value = "  +1.33%  "

# Row 50
$ws.Range("E50").Value = "  +1.30%  "

# Row 51
$ws.Range("E51").Value = "  +0.70%  "
